# Regenerate save_data column G (K = strikeouts) with recomputed values.
# This replaces the previous "Strike#" derived values with the new K values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 2
    6  = 0
    7  = 2
    8  = 3
    9  = 3
    10 = 0
    11 = 1
    12 = 2
    13 = 2
    14 = 1
    15 = 1
    17 = 2
    18 = 3
    19 = 0
    20 = 2
    21 = 2
    22 = 3
    23 = 2
    24 = 0
    26 = 1
    28 = 1
    29 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
